# Weekly update: insert this week's new Albahaca price record at the top
# of the data block (row 27), pushing all the existing records down by
# one row (old row 27 -> 28, ..., old row 65 -> 66).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows (27-65) down by inserting a new row at 27.
$ws.Rows.Item(27).Insert()

# Populate the newly inserted row 27 with this week's record.
$ws.Range("A27").Value = 8
$ws.Range("B27").Value = "Terminal La Palmera de La Serena"
$ws.Range("C27").Value = "Coquimbo"
$ws.Range("D27").Value = 44533
$ws.Range("E27").Value = 4
$ws.Range("F27").Value = 100112052
$ws.Range("G27").Value = "Albahaca"
$ws.Range("H27").Value = "Sin especificar"
$ws.Range("I27").Value = "Primera"
$ws.Range("J27").Value = 900
$ws.Range("K27").Value = 3000
$ws.Range("L27").Value = 4000
$ws.Range("M27").Value = 3500
$ws.Range("N27").Value = "$/paquete"
$ws.Range("O27").Value = "Región de Arica y Parinacota"
$ws.Range("P27").Value = 3500
$ws.Range("Q27").Value = 1
$ws.Range("R27").Value = "Hortaliza"
